$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_log")
# "time_log" is the workbook's active tab; $wb.ActiveSheet resolves the same sheet.

$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = Get-Date -Year 2023 -Month 6 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("B11").Value = "run_me clean, debug, & run"
$ws.Range("C11").Value = "Resolving discrepancies between raw and Calculations files cytotoxicity data"
$ws.Range("D11").Value = 2.42

$ws.Range("D12").Select()
